# Generate Report for handback
# Updates the localization-status workbook to reflect that the two
# handoff packages (57b8156c... and 92b3dd88...) have now been handed
# back and are in sync with en-US: the Status column moves from
# "Ready for handoff" to "Handed back: in sync with en-US", the
# "Latest Target File" / "Latest Handback File" columns (E/F) get
# populated (same files, since nothing changed), and the
# "Latest Handback DateTime" column (G) gets a real timestamp instead
# of the epoch placeholder.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: Status columns (B/C) for the two tracked files ---
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet ---
# Row 2: 57b8156c-cf9f-4f67-a157-8cfdef53e762
$wsZhCn.Range("B2").Value = $newStatus
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/ce844559eea2dcd244ca4bdf64d276d5479944c5/e2e/57b8156c-cf9f-4f67-a157-8cfdef53e762.md",
    "",
    "",
    "57b8156c-cf9f-4f67-a157-8cfdef53e762.md"
) | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ba7701aaa9dfa03f447b0e3cf1b8e05fe6fe3dbe/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/57b8156c-cf9f-4f67-a157-8cfdef53e762.f66315874c1f42410fcb52d24bdccd96ead7e29a.zh-cn.xlf",
    "",
    "",
    "57b8156c-cf9f-4f67-a157-8cfdef53e762.f66315874c1f42410fcb52d24bdccd96ead7e29a.zh-cn.xlf"
) | Out-Null
$wsZhCn.Range("G2").Value = "2016-01-11 07:53:52"

# Row 3: 92b3dd88-06b7-4bbf-acd7-ac7f81d3b112
$wsZhCn.Range("B3").Value = $newStatus
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/ce844559eea2dcd244ca4bdf64d276d5479944c5/e2e/92b3dd88-06b7-4bbf-acd7-ac7f81d3b112.md",
    "",
    "",
    "92b3dd88-06b7-4bbf-acd7-ac7f81d3b112.md"
) | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ba7701aaa9dfa03f447b0e3cf1b8e05fe6fe3dbe/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/92b3dd88-06b7-4bbf-acd7-ac7f81d3b112.ffb6f841966544fb26d211805f267cd32d2f57d5.zh-cn.xlf",
    "",
    "",
    "92b3dd88-06b7-4bbf-acd7-ac7f81d3b112.ffb6f841966544fb26d211805f267cd32d2f57d5.zh-cn.xlf"
) | Out-Null
$wsZhCn.Range("G3").Value = "2016-01-11 07:53:52"

# --- de-de sheet ---
# Row 2: 57b8156c-cf9f-4f67-a157-8cfdef53e762
$wsDeDe.Range("B2").Value = $newStatus
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/ce844559eea2dcd244ca4bdf64d276d5479944c5/e2e/57b8156c-cf9f-4f67-a157-8cfdef53e762.md",
    "",
    "",
    "57b8156c-cf9f-4f67-a157-8cfdef53e762.md"
) | Out-Null
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/702eb86955c90dbcd8b153329fc5bc1be3a5883d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/57b8156c-cf9f-4f67-a157-8cfdef53e762.f66315874c1f42410fcb52d24bdccd96ead7e29a.de-de.xlf",
    "",
    "",
    "57b8156c-cf9f-4f67-a157-8cfdef53e762.f66315874c1f42410fcb52d24bdccd96ead7e29a.de-de.xlf"
) | Out-Null
$wsDeDe.Range("G2").Value = "2016-01-11 07:54:26"

# Row 3: 92b3dd88-06b7-4bbf-acd7-ac7f81d3b112
$wsDeDe.Range("B3").Value = $newStatus
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/ce844559eea2dcd244ca4bdf64d276d5479944c5/e2e/92b3dd88-06b7-4bbf-acd7-ac7f81d3b112.md",
    "",
    "",
    "92b3dd88-06b7-4bbf-acd7-ac7f81d3b112.md"
) | Out-Null
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/702eb86955c90dbcd8b153329fc5bc1be3a5883d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/92b3dd88-06b7-4bbf-acd7-ac7f81d3b112.ffb6f841966544fb26d211805f267cd32d2f57d5.de-de.xlf",
    "",
    "",
    "92b3dd88-06b7-4bbf-acd7-ac7f81d3b112.ffb6f841966544fb26d211805f267cd32d2f57d5.de-de.xlf"
) | Out-Null
$wsDeDe.Range("G3").Value = "2016-01-11 07:54:26"
